$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column F (shifts old F:I -> G:J, preserving their content/format)
$ws.Columns("F").Insert()

# New column F matches the width of its neighbour column E
$ws.Columns("E:F").ColumnWidth = 20

# Header text for the newly inserted column F1 ("Detailed address (per corporate registry)")
$ws.Range("F1").Value = "상세주소`n(법인등기부등본상)"
$ws.Range("F1").WrapText = $true
$ws.Rows(1).RowHeight = 34.8

# Move the active selection to F2 (matches the post-edit cursor position)
$ws.Range("F2").Select() | Out-Null
